$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mai")

# Add the missing clock-in / clock-out times for row 33 (entry that was
# previously left blank) and mark the day's activity as "Coding".
$ws.Range("D33").Value = 0.80208333333333337
$ws.Range("E33").Value = 0.84375
$ws.Range("O33").Value = "Coding"

# Leave the cursor where the author ended up after filling the row in.
$ws.Range("G29").Select()
